# Updates crypto price/volume data rows 2-51 to reflect the latest scrape.
# (GitHub Actions scheduled refresh of cryptos.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.057.99"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.664.00"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'310.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.3906"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'0.3876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "'51.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("D10").Value = "'1.373"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "'1.000"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "'0.08510"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'24.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("D14").Value = "'7.238"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "'8.039"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.45%  "
$ws.Range("D16").Value = "'0.00001319"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "1.660.29"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "'94.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").Value = "'0.06999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").Value = "'20.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'7.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'13.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "24.045.10"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.511"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'3.153"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.67%  "
$ws.Range("D27").Value = "'22.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "'153.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").Value = "'141.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "'5.336"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "'7.842"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "1.842.72"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "'1.059"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.18%  "
$ws.Range("D35").Value = "'0.08189"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +4.73%  "
$ws.Range("D37").Value = "'11.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.74%  "
$ws.Range("D38").Value = "'6.735"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "'0.2720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("D40").Value = "'0.09168"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'13.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.43%  "
$ws.Range("D42").Value = "'0.7624"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").Value = "'1.431"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "'16.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("D45").Value = "'0.7048"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("D46").Value = "'2.515"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "'0.08333"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'135.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").Value = "'1.244"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.00%  "
